# Update "想去人数" (interested-count) figures on the "展览" and "全部类型"
# sheets to reflect the latest scrape, per the upstream data refresh.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 17
$ws1.Range("F3").Value = 1413
$ws1.Range("F7").Value = 11910
$ws1.Range("F8").Value = 4443
$ws1.Range("F11").Value = 28
$ws1.Range("F12").Value = 22
$ws1.Range("F13").Value = 2569
$ws1.Range("F15").Value = 166
$ws1.Range("F17").Value = 5164
$ws1.Range("F21").Value = 11390
$ws1.Range("F22").Value = 11383

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 17
$ws4.Range("F3").Value = 1413
$ws4.Range("F7").Value = 11910
$ws4.Range("F8").Value = 4443
$ws4.Range("F11").Value = 28
$ws4.Range("F12").Value = 22
$ws4.Range("F13").Value = 2569
$ws4.Range("F16").Value = 166
$ws4.Range("F18").Value = 5164
$ws4.Range("F22").Value = 11390
$ws4.Range("F23").Value = 11383
